$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace old evalsr-doc / evalsr-ext namespace prefixes and element names
# with the new phisr-doc / phisr-ext equivalents across the NIEM XPath column.
$map = @{
    "C2"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonName/nc30:PersonGivenName"
    "C3"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonName/nc30:PersonMiddleName"
    "C4"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonName/nc30:PersonSurName"
    "C5"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonBirthDate/nc30:Date"
    "C6"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/jxdm51:PersonRaceCode"
    "C7"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/jxdm51:PersonEthnicityCode"
    "C8"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/jxdm51:PersonSexCode"
    "C9"  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/phisr-ext:PersonTemporaryIdentification/nc30:IdentificationID"
    "C10" = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonName/nc30:PersonGivenName"
    "C11" = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonName/nc30:PersonMiddleName"
    "C12" = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonName/nc30:PersonSurName"
    "C13" = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonBirthDate/nc30:Date"
    "C14" = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/jxdm51:PersonSexCode"
}

foreach ($addr in $map.Keys) {
    $ws.Range($addr).Value = $map[$addr]
}

# Update the active selection to reflect where editing left the cursor.
$ws.Range("C16").Select()
